$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-300) holds a date serial number that was bumped by one day
# (45179 -> 45180, i.e. 2023-09-10 -> 2023-09-11) for every data row.
$ws.Range("C2:C300").Value = 45180
